# Insert two new rows at the top of the "Femacal de La Calera - Alcachofa" data
# block (rows 515:516), which pushes the existing rows 515:583 down to 517:585.
# Then populate the two new rows with the latest week's price entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before the current row 515 (shifts 515:583 -> 517:585)
$ws.Rows("515:516").Insert()

# --- New row 515: Alcachofa, Argentina(o), Primera ---
$ws.Cells.Item(515, 1).Value  = 3
$ws.Cells.Item(515, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(515, 3).Value  = "Coquimbo"
$ws.Cells.Item(515, 4).Value  = 45142
$ws.Cells.Item(515, 5).Value  = 5
$ws.Cells.Item(515, 6).Value  = 100112013
$ws.Cells.Item(515, 7).Value  = "Alcachofa"
$ws.Cells.Item(515, 8).Value  = "Argentina(o)"
$ws.Cells.Item(515, 9).Value  = "Primera"
$ws.Cells.Item(515, 10).Value = 60
$ws.Cells.Item(515, 11).Value = 13000
$ws.Cells.Item(515, 12).Value = 13000
$ws.Cells.Item(515, 13).Value = 13000
$ws.Cells.Item(515, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(515, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(515, 16).Value = 260
$ws.Cells.Item(515, 17).Value = 50
$ws.Cells.Item(515, 18).Value = "Hortaliza"

# --- New row 516: Alcachofa, Española, Primera ---
$ws.Cells.Item(516, 1).Value  = 3
$ws.Cells.Item(516, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(516, 3).Value  = "Coquimbo"
$ws.Cells.Item(516, 4).Value  = 45142
$ws.Cells.Item(516, 5).Value  = 5
$ws.Cells.Item(516, 6).Value  = 100112013
$ws.Cells.Item(516, 7).Value  = "Alcachofa"
$ws.Cells.Item(516, 8).Value  = "Española"
$ws.Cells.Item(516, 9).Value  = "Primera"
$ws.Cells.Item(516, 10).Value = 55
$ws.Cells.Item(516, 11).Value = 14000
$ws.Cells.Item(516, 12).Value = 14000
$ws.Cells.Item(516, 13).Value = 14000
$ws.Cells.Item(516, 14).Value = "`$/caja 30 unidades"
$ws.Cells.Item(516, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(516, 16).Value = 467
$ws.Cells.Item(516, 17).Value = 30
$ws.Cells.Item(516, 18).Value = "Hortaliza"
